$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so values such as
# "193.32" are not auto-converted into numeric cells by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.746.60"
$ws.Range("E2").Value = "  +3.35%  "
$ws.Range("D3").Value = "3.382.72"
$ws.Range("E3").Value = "  +4.82%  "
$ws.Range("D5").Value = "193.32"
$ws.Range("E5").Value = "  +6.17%  "
$ws.Range("D6").Value = "594.46"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.609"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("D9").Value = "0.135"
$ws.Range("E9").Value = "  +3.56%  "
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("D11").Value = "0.424"
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("D12").Value = "3.970.59"
$ws.Range("E12").Value = "  +4.82%  "
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "28.79"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").Value = "69.742.85"
$ws.Range("E15").Value = "  +3.25%  "
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").Value = "3.379.62"
$ws.Range("E17").Value = "  +6.08%  "
$ws.Range("D18").Value = "451.31"
$ws.Range("E18").Value = "  +14.51%  "
$ws.Range("D19").Value = "5.86"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "13.80"
$ws.Range("E20").Value = "  +2.78%  "
$ws.Range("D21").Value = "7.85"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("D22").Value = "73.54"
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "3.513.27"
$ws.Range("E24").Value = "  +4.36%  "
$ws.Range("D26").Value = "0.520"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("E27").Value = "  +4.63%  "
$ws.Range("D28").Value = "9.63"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  +2.85%  "
$ws.Range("D31").Value = "23.28"
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("D32").Value = "5.66"
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").Value = "1.30"
$ws.Range("E33").Value = "  +4.08%  "
$ws.Range("D34").Value = "7.06"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "1.53"
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("D37").Value = "164.83"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("D39").Value = "27.31"
$ws.Range("E39").Value = "  +3.86%  "
$ws.Range("D40").Value = "0.825"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").Value = "4.62"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").Value = "6.53"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").Value = "2.748.69"
$ws.Range("E43").Value = "  +5.46%  "
$ws.Range("D44").Value = "2.56"
$ws.Range("E44").Value = "  +4.10%  "
$ws.Range("D45").Value = "25.69"
$ws.Range("E45").Value = "  +4.35%  "
$ws.Range("D46").Value = "0.0692"
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("D47").Value = "345.91"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("D48").Value = "40.85"
$ws.Range("D49").Value = "0.0287"
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("E50").Value = "  +7.65%  "
$ws.Range("E51").Value = "  +5.03%  "